# Weekly fruit/vegetable price update: insert two new observation rows for
# "Ají" (Vega Central Mapocho de Santiago) right after the existing row 166,
# shifting the previous rows 167-189 down to 169-191, and populate the two
# new rows with their data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at position 167 (pushes old row167.. down by 2)
$ws.Rows.Item(167).Insert()
$ws.Rows.Item(167).Insert()

# --- New row 167 ---------------------------------------------------------
$ws.Cells.Item(167, 1).Value = 9
$ws.Cells.Item(167, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(167, 3).Value = "Metropolitana"
$ws.Cells.Item(167, 4).Value = 44522
$ws.Cells.Item(167, 5).Value = 13
$ws.Cells.Item(167, 6).Value = 100112021
$ws.Cells.Item(167, 7).Value = "Ají"
$ws.Cells.Item(167, 8).Value = "Inferno"
$ws.Cells.Item(167, 9).Value = "Primera"
$ws.Cells.Item(167, 10).Value = 52
$ws.Cells.Item(167, 11).Value = 18000
$ws.Cells.Item(167, 12).Value = 20000
$ws.Cells.Item(167, 13).Value = 19000
$ws.Cells.Item(167, 14).Value = "$/caja 12 kilos"
$ws.Cells.Item(167, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(167, 16).Value = 1583
$ws.Cells.Item(167, 17).Value = 12
$ws.Cells.Item(167, 18).Value = "Hortaliza"

# --- New row 168 ---------------------------------------------------------
$ws.Cells.Item(168, 1).Value = 9
$ws.Cells.Item(168, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(168, 3).Value = "Metropolitana"
$ws.Cells.Item(168, 4).Value = 44522
$ws.Cells.Item(168, 5).Value = 13
$ws.Cells.Item(168, 6).Value = 100112021
$ws.Cells.Item(168, 7).Value = "Ají"
$ws.Cells.Item(168, 8).Value = "Inferno"
$ws.Cells.Item(168, 9).Value = "Segunda"
$ws.Cells.Item(168, 10).Value = 25
$ws.Cells.Item(168, 11).Value = 15000
$ws.Cells.Item(168, 12).Value = 15000
$ws.Cells.Item(168, 13).Value = 15000
$ws.Cells.Item(168, 14).Value = "$/caja 12 kilos"
$ws.Cells.Item(168, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(168, 16).Value = 1250
$ws.Cells.Item(168, 17).Value = 12
$ws.Cells.Item(168, 18).Value = "Hortaliza"
